$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) date column C for existing data rows
#    2..77 from 2023-09-21 (45190) to 2023-09-23 (45192).
$ws.Range("C2:C77").Value = 45192

# 2. Row 77 picks up an explicit row height in the saved file — set it
#    explicitly so it is written out with ht="15" customHeight="1".
$ws.Rows.Item(77).RowHeight = 15

# 3. Append the new record as row 78.
$ws.Range("A78").Value = "A 44496-2023"

$ws.Range("B78").Value = 45189
$ws.Range("B78").NumberFormat = "YYYY-MM-DD"

$ws.Range("C78").Value = 45192
$ws.Range("C78").NumberFormat = "YYYY-MM-DD"

$ws.Range("D78").Value = "SKÅNE LÄN"
$ws.Range("E78").Value = "ESLÖV"

$ws.Range("G78").Value = 0.5
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("O78").Value = 0
$ws.Range("P78").Value = 0
$ws.Range("Q78").Value = 0

# R column uses a wrap-text style across the sheet (even when empty).
$ws.Range("R78").WrapText = $true
